$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Jones"
$ws.Range("B3").Value = "2 Adults, 3 Children"
$ws.Range("C3").Value = "Peaches"
$ws.Range("D3").Value = "Peanuts"
$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "Yes"
$ws.Range("G3").Value = "Yes"
$ws.Range("H3").Value = "Yes"
$ws.Range("I3").Value = "Yes"
$ws.Range("J3").Value = "No"
$ws.Range("K3").Value = "Yes"
$ws.Range("L3").Value = "Yes"

$ws.Range("C3").Select() | Out-Null
